$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Steel Annealed (column C) values changed for the 20C and 0C rows
$ws.Range("C4").Value = 0
$ws.Range("C5").Value = 40

# Lowest temperature data point updated to the (approx) real value, -195C
$ws.Range("A9").Value = -195

# Move / leave the active selection where the author left off editing
$ws.Range("E12").Select()
